# Insert a new weekly data row above row 118 (Jengibre / Vega Modelo de Temuco),
# pushing the existing rows 118-126 down to 119-127, and populate the new
# row 118 with the latest price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 118 - shifts rows 118:126 down to 119:127
$ws.Rows.Item(118).Insert()

# Fill in the new row 118 with the latest weekly record
$ws.Cells.Item(118, 1).Value = 10
$ws.Cells.Item(118, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(118, 3).Value = "La Araucanía"
$ws.Cells.Item(118, 4).Value = 44578
$ws.Cells.Item(118, 5).Value = 9
$ws.Cells.Item(118, 6).Value = 100114007
$ws.Cells.Item(118, 7).Value = "Jengibre"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 200
$ws.Cells.Item(118, 11).Value = 20000
$ws.Cells.Item(118, 12).Value = 20000
$ws.Cells.Item(118, 13).Value = 20000
$ws.Cells.Item(118, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(118, 15).Value = "Perú"
$ws.Cells.Item(118, 16).Value = 1538
$ws.Cells.Item(118, 17).Value = 13
$ws.Cells.Item(118, 18).Value = "Hortaliza"
